$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=103; Date=45236; Time="20:39"; Price=2699;    Site="amazon";        Cor="preto" },
    @{ Row=104; Date=45236; Time="20:40"; Price=2625;    Site="mercado livre"; Cor="preto" },
    @{ Row=105; Date=45237; Time="20:27"; Price=2607.01; Site="amazon";        Cor="preto" },
    @{ Row=106; Date=45237; Time="20:27"; Price=2399;    Site="mercado livre"; Cor="preto" }
)

$dateFmt = $ws.Cells.Item(102, 1).NumberFormat

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 1).NumberFormat = $dateFmt
    $ws.Cells.Item($r.Row, 2).Value = $r.Time
    $ws.Cells.Item($r.Row, 3).Value = $r.Price
    $ws.Cells.Item($r.Row, 4).Value = $r.Site
    $ws.Cells.Item($r.Row, 5).Value = $r.Cor
}
